{"js": "const tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\n// Row 0 is the header (\"S.N.\", \"Roll No.\", ...); row 1 is the single data\n// row (S.N. 1, roll no. 077bct098, ...) whose values need to be cleared\n// out, leaving the empty runs (and the _GoBack bookmark) in place.\nconst row = rows.items[1];\nconst cells = row.cells;\ncells.load(\"items\");\nawait context.sync();\n\nfor (const cell of cells.items) {\n  // Replace the cell's text in place (via its body range) instead of\n  // TableCell.value, which rewrites the whole paragraph and would drop\n  // the paraId/rsid attributes and the bookmark in the first cell.\n  const range = cell.body.getRange();\n  range.insertText(\"\", \"Replace\");\n}\nawait context.sync();\n", "ps1": "$doc = $word.ActiveDocument\n\n# The document has a single results table. Row 1 is the bold header\n# (\"S.N.\", \"Roll No.\", \"Exam Roll No.\", \"Name\", ...); row 2 is the lone\n# data row (S.N. 1 / 077bct098 / sdf56 / darpan / 30 / 70 / 100) whose\n# values need to be wiped out, leaving the empty runs (and the _GoBack\n# bookmark in the first cell) in place.\n$table = $doc.Tables.Item(1)\n$row = $table.Rows.Item(2)\n\nfor ($i = 1; $i -le $row.Cells.Count; $i++) {\n    $cell = $row.Cells.Item($i)\n    $cell.Range.Text = \"\"\n}\n"}
